# Add data for 2021-11-24: update the "through" date in the sheet name and
# the November row label, and refresh the 2016/2017/2018/2020/2021 monthly
# (through 11-16) and Total rows with the new arrest/no_arrest counts and
# recomputed arrest rates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab to reflect the new "through" date.
$ws.Name = "Through 2021-11-16"

# Update the row-13 label.
$ws.Range("A13").Value = "November (through 11-16)"

# Row 13 ("November (through 11-16)") — no_arrest_made + arrest_rate per year.
$ws.Range("F13").Value = 36
$ws.Range("G13").Value = 0.0769

$ws.Range("I13").Value = 66
$ws.Range("J13").Value = 0.0294

$ws.Range("L13").Value = 29
$ws.Range("M13").Value = 0.1471

$ws.Range("R13").Value = 90
$ws.Range("S13").Value = 0.0323

$ws.Range("U13").Value = 113
$ws.Range("V13").Value = 0.0088

# Row 14 ("Total") — no_arrest_made + arrest_rate per year.
$ws.Range("F14").Value = 470
$ws.Range("G14").Value = 0.1048

$ws.Range("I14").Value = 715
$ws.Range("J14").Value = 0.081

$ws.Range("L14").Value = 578
$ws.Range("M14").Value = 0.1094

$ws.Range("R14").Value = 1093
$ws.Range("S14").Value = 0.0496

$ws.Range("U14").Value = 1467
$ws.Range("V14").Value = 0.0572
